# Apply the "ExSim dashboards" data refresh to the Logistics Dashboard workbook.
# Zone capacities / opening inventory / rent figures are revised on
# ROUTE_CONFIG and INVENTORY_TETRIS, the shipment-cost benchmark values on
# SHIPMENT_BUILDER are refreshed (including the conditional-format rule),
# and two stale "Last Period Rent" labels (zones with 0 capacity) are
# removed entirely.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ROUTE_CONFIG — zone capacity table (B13:B15)
# ---------------------------------------------------------------------
$routeConfig = $wb.Worksheets.Item("ROUTE_CONFIG")
$routeConfig.Range("B13").Value = 4800   # Center
$routeConfig.Range("B14").Value = 2500   # West
$routeConfig.Range("B15").Value = 2000   # North

# ---------------------------------------------------------------------
# INVENTORY_TETRIS — per-zone headers / opening inventory / capacity
# ---------------------------------------------------------------------
$tetris = $wb.Worksheets.Item("INVENTORY_TETRIS")

# Center zone (row 4-5)
$tetris.Range("A4").Value = "═══ CENTER ZONE (Capacity: 4,800) ═══"
$tetris.Range("I4").Value = "Last Period Rent: $38,400"
$tetris.Range("B5").Value = 3410
$tetris.Range("E5").Value = 4800

# West zone (row 18-19)
$tetris.Range("A18").Value = "═══ WEST ZONE (Capacity: 2,500) ═══"
$tetris.Range("I18").Value = "Last Period Rent: $20,000"
$tetris.Range("B19").Value = 1886
$tetris.Range("E19").Value = 2500

# North zone (row 32-33)
$tetris.Range("A32").Value = "═══ NORTH ZONE (Capacity: 2,000) ═══"
$tetris.Range("I32").Value = "Last Period Rent: $16,000"
$tetris.Range("B33").Value = 1720
$tetris.Range("E33").Value = 2000

# East / South zones keep 0 capacity, but their stale "Last Period Rent"
# labels are removed outright (the whole cell disappears, not just its
# text) since those zones have no rent to report.
$tetris.Range("I46").Clear()
$tetris.Range("I60").Clear()

# ---------------------------------------------------------------------
# SHIPMENT_BUILDER — route cost-per-unit benchmarks + red-flag threshold
# ---------------------------------------------------------------------
$shipment = $wb.Worksheets.Item("SHIPMENT_BUILDER")
$shipment.Range("N6").Value = 13.2
$shipment.Range("N7").Value = 16.31992540096979
$shipment.Range("A27").Value = "⚠️ Red cells = Cost >20% above avg benchmark ($14.76 avg)"

$costRange = $shipment.Range("F6:F25")
$costRule = $costRange.FormatConditions.Item(1)
$costRule.Formula1 = "=F6>17.711955240581872"
